$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.330193400382996
$ws.Range("B1").Value = 1.665873289108276
$ws.Range("C1").Value = 4.109106540679932
$ws.Range("D1").Value = 3.177010297775269
$ws.Range("E1").Value = 1.111135125160217
